$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.138377666473389
$ws.Range("B1").Value = 2.415556192398071
$ws.Range("C1").Value = 5.19922924041748
$ws.Range("D1").Value = 2.223955154418945
$ws.Range("E1").Value = 1.2546706199646
